$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 64
$ws.Range("I9").Value = 58.333332
$ws.Range("J9").Value = 68.25
$ws.Range("K9").Value = 58.333332
$ws.Range("L9").Value = 68.25
$ws.Range("M9").Value = 110.666668
$ws.Range("N9").Value = -406.25
$ws.Range("H40").Value = 2416.6667
$ws.Range("J40").Value = 2500
$ws.Range("L40").Value = 2500
$ws.Range("N40").Value = -2850
$ws.Range("H43").Value = 8193.5
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H70").Value = 17500.5
$ws.Range("J70").Value = 17500.5
$ws.Range("L70").Value = 52501.5
$ws.Range("N70").Value = -53041.5
$ws.Range("H73").Value = 17500.5
$ws.Range("J73").Value = 17500.5
$ws.Range("L73").Value = 52501.5
$ws.Range("N73").Value = -54373.5
$ws.Range("H76").Value = 5365.826
$ws.Range("I76").Value = 4246.8184
$ws.Range("J76").Value = 6391.5835
$ws.Range("K76").Value = 4246.8184
$ws.Range("L76").Value = 6391.5835
$ws.Range("M76").Value = -3931.8184
$ws.Range("N76").Value = -7021.5835
$ws.Range("H79").Value = 5365.826
$ws.Range("I79").Value = 4246.8184
$ws.Range("J79").Value = 6391.5835
$ws.Range("K79").Value = 4246.8184
$ws.Range("L79").Value = 6391.5835
$ws.Range("M79").Value = -3154.8184
$ws.Range("N79").Value = -8575.583500000001
$ws.Range("H80").Value = 510.36365
$ws.Range("I80").Value = 383
$ws.Range("K80").Value = 1149
$ws.Range("M80").Value = -151
$ws.Range("H83").Value = 510.36365
$ws.Range("I83").Value = 383
$ws.Range("K83").Value = 3447
$ws.Range("M83").Value = 1545
$ws.Range("H116").Value = 6155.154
$ws.Range("J116").Value = 6971.3335
$ws.Range("L116").Value = 6971.3335
$ws.Range("N116").Value = -13855.3335
$ws.Range("H132").Value = 1067.1471
$ws.Range("I132").Value = 1141.6774
$ws.Range("J132").Value = 297
$ws.Range("K132").Value = 3425.0322
$ws.Range("L132").Value = 891
$ws.Range("M132").Value = -895.0322000000001
$ws.Range("N132").Value = -5951
$ws.Range("H138").Value = 4102.1333
$ws.Range("I138").Value = 1864.875
$ws.Range("J138").Value = 4585.8647
$ws.Range("K138").Value = 5594.625
$ws.Range("L138").Value = 13757.5941
$ws.Range("M138").Value = -454.625
$ws.Range("N138").Value = -24037.5941

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1994.2858
$ws.Range("I2").Value = 1493.3334
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 1493.3334
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -1380.3334
$ws.Range("N2").Value = -5226
$ws.Range("H31").Value = 3850
$ws.Range("I31").Value = 3850
$ws.Range("K31").Value = 3850
$ws.Range("M31").Value = -3556
$ws.Range("H45").Value = 1686.875
$ws.Range("I45").Value = 1675.5238
$ws.Range("J45").Value = 1766.3334
$ws.Range("K45").Value = 1675.5238
$ws.Range("L45").Value = 1766.3334
$ws.Range("M45").Value = -1298.5238
$ws.Range("N45").Value = -2520.3334
$ws.Range("H116").Value = 1994.2858
$ws.Range("I116").Value = 1493.3334
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 1493.3334
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = 800.6666
$ws.Range("N116").Value = -9588
$ws.Range("H122").Value = 437541.12
$ws.Range("I122").Value = 502572.3
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 1507716.9
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -1505266.9
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 3022.3076
$ws.Range("I132").Value = 2662.7273
$ws.Range("K132").Value = 7988.1819
$ws.Range("M132").Value = -5458.1819

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1994.2858
$ws.Range("I3").Value = 1493.3334
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 1493.3334
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -1379.3334
$ws.Range("N3").Value = -5228
$ws.Range("H80").Value = 327.84616
$ws.Range("I80").Value = 379.66666
$ws.Range("J80").Value = 211.25
$ws.Range("K80").Value = 379.66666
$ws.Range("L80").Value = 211.25
$ws.Range("M80").Value = 618.33334
$ws.Range("N80").Value = -2207.25
$ws.Range("H83").Value = 327.84616
$ws.Range("I83").Value = 379.66666
$ws.Range("J83").Value = 211.25
$ws.Range("K83").Value = 1898.3333
$ws.Range("L83").Value = 1056.25
$ws.Range("M83").Value = 3093.6667
$ws.Range("N83").Value = -11040.25
$ws.Range("H99").Value = 3361.4783
$ws.Range("J99").Value = 3744.1667
$ws.Range("L99").Value = 3744.1667
$ws.Range("N99").Value = -6740.1667

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 871.8182
$ws.Range("I16").Value = 854.44446
$ws.Range("J16").Value = 950
$ws.Range("K16").Value = 854.44446
$ws.Range("L16").Value = 950
$ws.Range("M16").Value = -567.44446
$ws.Range("N16").Value = -1524
$ws.Range("H74").Value = 42226.6
$ws.Range("I74").Value = 40273.332
$ws.Range("J74").Value = 45156.5
$ws.Range("K74").Value = 40273.332
$ws.Range("L74").Value = 45156.5
$ws.Range("M74").Value = -39399.332
$ws.Range("N74").Value = -46904.5
$ws.Range("H77").Value = 42226.6
$ws.Range("I77").Value = 40273.332
$ws.Range("J77").Value = 45156.5
$ws.Range("K77").Value = 120819.996
$ws.Range("L77").Value = 135469.5
$ws.Range("M77").Value = -116451.996
$ws.Range("N77").Value = -144205.5
$ws.Range("H99").Value = 13234.969
$ws.Range("I99").Value = 8920.385
$ws.Range("K99").Value = 8920.385
$ws.Range("M99").Value = -7422.385
$ws.Range("H105").Value = 1799.2
$ws.Range("I105").Value = 784.7143
$ws.Range("K105").Value = 784.7143
$ws.Range("M105").Value = 962.2857
$ws.Range("H113").Value = 871.8182
$ws.Range("I113").Value = 854.44446
$ws.Range("J113").Value = 950
$ws.Range("K113").Value = 854.44446
$ws.Range("L113").Value = 950
$ws.Range("M113").Value = 1315.55554
$ws.Range("N113").Value = -5290
$ws.Range("H126").Value = 13234.969
$ws.Range("I126").Value = 8920.385
$ws.Range("K126").Value = 26761.155
$ws.Range("M126").Value = -24291.155
$ws.Range("H132").Value = 2499.3
$ws.Range("I132").Value = 1370.75
$ws.Range("K132").Value = 4112.25
$ws.Range("M132").Value = -1582.25
$ws.Range("H134").Value = 2831.9546
$ws.Range("I134").Value = 2349.389
$ws.Range("K134").Value = 7048.167
$ws.Range("M134").Value = -4513.167

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1334260.5
$ws.Range("J4").Value = 1591.5714
$ws.Range("L4").Value = 4774.7142
$ws.Range("N4").Value = -4998.7142
$ws.Range("H17").Value = 93.22221999999999
$ws.Range("J17").Value = 99.333336
$ws.Range("L17").Value = 298.000008
$ws.Range("N17").Value = -636.000008
$ws.Range("H92").Value = 999
$ws.Range("J92").Value = 999
$ws.Range("L92").Value = 2997
$ws.Range("N92").Value = -5493
$ws.Range("H106").Value = 9999.799999999999
$ws.Range("J106").Value = 9999.799999999999
$ws.Range("L106").Value = 29999.4
$ws.Range("N106").Value = -31891.4
$ws.Range("H140").Value = 2247
$ws.Range("I140").Value = 2591.1
$ws.Range("J140").Value = 1100
$ws.Range("K140").Value = 7773.299999999999
$ws.Range("L140").Value = 3300
$ws.Range("M140").Value = -2593.299999999999
$ws.Range("N140").Value = -13660

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 10010002
$ws.Range("I12").Value = 20000000
$ws.Range("J12").Value = 20003.5
$ws.Range("K12").Value = 20000000
$ws.Range("L12").Value = 20003.5
$ws.Range("M12").Value = -19999860
$ws.Range("N12").Value = -20283.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7324.4
$ws.Range("I122").Value = 7655.75
$ws.Range("J122").Value = 5999
$ws.Range("K122").Value = 22967.25
$ws.Range("L122").Value = 17997
$ws.Range("M122").Value = -20517.25
$ws.Range("N122").Value = -22897

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 630.7857
$ws.Range("I107").Value = 416.14285
$ws.Range("K107").Value = 1248.42855
$ws.Range("M107").Value = 671.5714499999999
$ws.Range("H113").Value = 537.3125
$ws.Range("I113").Value = 660.4545000000001
$ws.Range("J113").Value = 266.4
$ws.Range("K113").Value = 1981.3635
$ws.Range("L113").Value = 799.1999999999999
$ws.Range("M113").Value = 188.6364999999998
$ws.Range("N113").Value = -5139.2
$ws.Range("H126").Value = 4714.8335
$ws.Range("I126").Value = 4356.8
$ws.Range("K126").Value = 13070.4
$ws.Range("M126").Value = -10600.4
